{"js": "// Replace the three-digit-divided-by-one-digit division problems with a\n// newly generated set of problems, preserving every other aspect of the\n// document (formatting, table structure, blank spacer rows, etc.).\n//\n// The document contains a single table; the five \"content\" rows (every\n// 4th row starting at row 0) each hold five division expressions, one per\n// cell. The remaining rows are blank spacer rows and are left untouched.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// New values, in row-major order, for the five content rows x five columns.\nconst newValues = [\n  [\"918\u00f74=229, 2\", \"374\u00f73=124, 2\", \"221\u00f79=24, 5\", \"793\u00f76=132, 1\", \"367\u00f72=183, 1\"],\n  [\"229\u00f73=76, 1\", \"600\u00f76=100, 0\", \"640\u00f76=106, 4\", \"560\u00f76=93, 2\", \"756\u00f78=94, 4\"],\n  [\"359\u00f79=39, 8\", \"830\u00f73=276, 2\", \"438\u00f73=146, 0\", \"670\u00f73=223, 1\", \"904\u00f77=129, 1\"],\n  [\"679\u00f75=135, 4\", \"352\u00f72=176, 0\", \"559\u00f75=111, 4\", \"635\u00f73=211, 2\", \"938\u00f78=117, 2\"],\n  [\"289\u00f77=41, 2\", \"502\u00f72=251, 0\", \"822\u00f79=91, 3\", \"584\u00f73=194, 2\", \"988\u00f74=247, 0\"],\n];\n\n// The content rows are located at table row indices 0, 4, 8, 12, 16.\nconst contentRowIndices = [0, 4, 8, 12, 16];\n\nfor (let r = 0; r < contentRowIndices.length; r++) {\n  const tableRowIndex = contentRowIndices[r];\n  for (let c = 0; c < newValues[r].length; c++) {\n    const cell = table.getCell(tableRowIndex, c);\n    cell.value = newValues[r][c];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the three-digit-divided-by-one-digit division problems with a\n# newly generated set of problems, preserving every other aspect of the\n# document (formatting, table structure, blank spacer rows, etc.).\n#\n# The document contains a single table; the five \"content\" rows (every\n# 4th row starting at row 1 in COM's 1-based indexing) each hold five\n# division expressions, one per cell. The remaining rows are blank spacer\n# rows and are left untouched.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# New values, in row-major order, for the five content rows x five columns.\n$newValues = @(\n    @(\"918\u00f74=229, 2\", \"374\u00f73=124, 2\", \"221\u00f79=24, 5\", \"793\u00f76=132, 1\", \"367\u00f72=183, 1\"),\n    @(\"229\u00f73=76, 1\", \"600\u00f76=100, 0\", \"640\u00f76=106, 4\", \"560\u00f76=93, 2\", \"756\u00f78=94, 4\"),\n    @(\"359\u00f79=39, 8\", \"830\u00f73=276, 2\", \"438\u00f73=146, 0\", \"670\u00f73=223, 1\", \"904\u00f77=129, 1\"),\n    @(\"679\u00f75=135, 4\", \"352\u00f72=176, 0\", \"559\u00f75=111, 4\", \"635\u00f73=211, 2\", \"938\u00f78=117, 2\"),\n    @(\"289\u00f77=41, 2\", \"502\u00f72=251, 0\", \"822\u00f79=91, 3\", \"584\u00f73=194, 2\", \"988\u00f74=247, 0\")\n)\n\n# The content rows are located at 1-based table row indices 1, 5, 9, 13, 17.\n$contentRowIndices = @(1, 5, 9, 13, 17)\n\nfor ($r = 0; $r -lt $contentRowIndices.Length; $r++) {\n    $tableRowIndex = $contentRowIndices[$r]\n    for ($c = 0; $c -lt $newValues[$r].Length; $c++) {\n        $cell = $t.Cell($tableRowIndex, $c + 1)\n        $cell.Range.Text = $newValues[$r][$c]\n    }\n}\n"}
